# Fruta / hortaliza, semanal
# Weekly refresh of the rolling price-history window (rows 16-24) for
# "Agrícola del Norte S.A. de Arica - Cilantro": drop the oldest entry
# (old row 16) and shift the remaining entries up by one row, then set
# the new date for the row that now represents the newest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) -- each row N takes on what
# used to be row N+1's data.
$rows = @(
    @{ Row = 16; D = 44218; J = 120; K = 1400; L = 1500; M = 1450; P = 725 },
    @{ Row = 17; D = 44417; J = 250; K = 1000; L = 1200; M = 1100; P = 550 },
    @{ Row = 18; D = 44295; J = 300; K = 1800; L = 2000; M = 1900; P = 950 },
    @{ Row = 19; D = 44175; J = 150; K = 1300; L = 1300; M = 1300; P = 650 },
    @{ Row = 20; D = 44355; J = 250; K = 2500; L = 3000; M = 2750; P = 1375 },
    @{ Row = 21; D = 44278; J = 200; K = 1000; L = 1200; M = 1100; P = 550 },
    @{ Row = 22; D = 44382; J = 300; K = 2400; L = 2500; M = 2450; P = 1225 },
    @{ Row = 23; D = 44284; J = 300; K = 900;  L = 1000; M = 950;  P = 475 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("P$row").Value = $r.P
}

# Row 24 gets the new (latest) date for the appended data point.
$ws.Range("D24").Value = 44482
